$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the footer note text (C9) to mention herramientas, tareas, empleados
$ws.Range("C9").Value = "Footer con últimos empleados, herramientas, tareas, empleados (extra)"

# Bump the date in D9 by one day (45773 -> 45774)
$ws.Range("D9").Value = 45774

# Move the active selection to D10 (matches the saved sheet view state)
$ws.Range("D10").Select()
